# novos graficos e fixes
# Rename the "Outro" category to "Outras" in the social-networks table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "Outras"
$ws.Range("A11").Select()
